# Updates the "cryptos" price list with refreshed price / volume(1h) figures,
# and swaps the Stacks / EnergySwap rows (49 and 50) to match the new ranking
# order, as produced by the scheduled GitHub Actions scraper run.
#
# Each target cell is written with a leading apostrophe so Excel stores the
# value as literal text (matching the original inlineStr cells) instead of
# silently re-interpreting numeric-looking strings (e.g. "1.00", "6.42") as
# numbers. The Style is then reset to "Normal" so no stray number-format /
# style index is left attached to the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.529.39"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.22%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.778.36"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.18%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D5").Value = "'597.65"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.28%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'164.54"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -1.33%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  -0.09%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  -1.20%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  -1.06%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  +0.40%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'6.42"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +1.29%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  -1.84%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'35.50"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -1.40%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'4.413.72"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -0.21%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'3.808.26"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -0.50%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'67.584.80"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.07%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'18.26"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -1.59%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = "'  +1.77%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'7.03"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -0.37%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'459.65"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +0.13%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'9.70"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -2.84%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  -0.51%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  -4.87%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'82.38"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -1.28%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'11.98"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -0.92%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'2.09"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -0.95%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  -0.07%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'9.94"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -0.52%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'3.927.73"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -0.26%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'7.44"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +3.27%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D32").Value = "'2.19"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -2.62%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'29.03"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -1.84%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'1.00"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -0.04%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'8.95"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -1.09%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'0.0986"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -1.15%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  +0.57%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'3.24"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -3.07%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.987"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -0.51%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  -0.30%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  -0.05%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D43").Value = "'43.49"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -0.68%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'47.36"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -1.27%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.297"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -0.10%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'151.72"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +0.98%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'8.31"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +0.60%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'1.35"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +8.46%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("B49").Value = "'Stacks"
$ws.Range("B49").Style = "Normal"
$ws.Range("C49").Value = "'https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("C49").Style = "Normal"
$ws.Range("D49").Value = "'1.84"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +1.62%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("B50").Value = "'EnergySwap"
$ws.Range("B50").Style = "Normal"
$ws.Range("C50").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("C50").Style = "Normal"
$ws.Range("D50").Value = "'26.89"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +0.19%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'391.52"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +0.75%  "
$ws.Range("E51").Style = "Normal"
